# Generate Report for Handback
# Updates the handback-status workbook with the latest handoff/handback
# generation timestamps for the file "1ae2545b-cd7f-4fe1-99b4-a33991af9af0.md"
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for row 2 (1ae2545b...md)
$wsOverview.Range("G2").Value = "2016-08-12 03:09:43"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsZhCn.Range("H2").Value = "2016-08-12 03:09:37"
$wsZhCn.Range("K2").Value = "2016-08-12 03:09:53"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsDeDe.Range("H2").Value = "2016-08-12 03:09:43"
$wsDeDe.Range("K2").Value = "2016-08-12 03:10:03"
